$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 385 (shifts existing rows 385..455 down to 386..456)
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new weekly record
$ws.Cells.Item(385, 1).Value = 5
$ws.Cells.Item(385, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(385, 3).Value = "Maule"
$ws.Cells.Item(385, 4).Value = 44637
$ws.Cells.Item(385, 5).Value = 7
$ws.Cells.Item(385, 6).Value = "Fruta"
$ws.Cells.Item(385, 7).Value = 100102
$ws.Cells.Item(385, 8).Value = "Cítricos"
$ws.Cells.Item(385, 9).Value = 100102005
$ws.Cells.Item(385, 10).Value = "Naranja"
$ws.Cells.Item(385, 11).Value = "Valencia"
$ws.Cells.Item(385, 12).Value = "Primera"
$ws.Cells.Item(385, 13).Value = 300
$ws.Cells.Item(385, 14).Value = 10000
$ws.Cells.Item(385, 15).Value = 10000
$ws.Cells.Item(385, 16).Value = 10000
$ws.Cells.Item(385, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(385, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(385, 19).Value = 667
$ws.Cells.Item(385, 20).Value = 15
